# Apply the "build area vs extrusion length" feature update to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update input values -------------------------------------------------
$ws.Range("A2").Value = 254
$ws.Range("D2").Value = 125
$ws.Range("D3").Value = 177
$ws.Range("A7").Value = 21.1

# --- Update the label in E14 to reflect the new DELTA_RADIUS naming -----
$text = "DELTA_RADIUS/ Horizontal radius: Diagonal Rod Horizontal length "
$boldLen = 31
$e14 = $ws.Range("E14")
$e14.Value = $text
$bold = $e14.Characters(1, $boldLen)
$bold.Font.Bold = $true
$bold.Font.Size = 10
$bold.Font.Name = "Arial"
$bold.Font.Color = 0
$rest = $e14.Characters($boldLen + 1, $text.Length - $boldLen)
$rest.Font.Bold = $false
$rest.Font.Size = 10
$rest.Font.Name = "Arial"
$rest.Font.Color = 0

# --- Move the active selection to A8 -------------------------------------
$ws.Activate()
$ws.Range("A8").Select()
